$d = $word.ActiveDocument

# Locate the paragraph that contains "LOB1012: Estatística (Requisito)".
# The three paragraphs that follow it need to be removed entirely:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. the "© 2020 ... Creative Commons Attribution" paragraph
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "LOB1012*Estat*stica*Requisito*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $p1 = $target.Next()
    $p2 = $p1.Next()
    $p3 = $p2.Next()

    # Delete from the last paragraph back to the first so ranges/indices
    # of the earlier paragraphs stay valid while deleting.
    $p3.Range.Delete()
    $p2.Range.Delete()
    $p1.Range.Delete()
}
